$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the remaining row 7 data (finishing the -30C data point)
$ws.Range("B7").Value = 8279
$ws.Range("C7").Value = 3.98548
$ws.Range("D7").Value = 3.13432
$ws.Range("E7").Value = -0.00104
$ws.Range("F7").Value = -4.49667

# G7 already holds the shared formula from G3:G7; now that E7/F7 are
# non-zero it will recalculate away from the #DIV/0! error.
$wb.Application.Calculate()

# Update the active selection to reflect where the author left off editing
$ws.Range("F8").Select()
